# Daily update at 8 AM UTC — append the next day's row to the Wins Over
# Time sheet (row 51: 2025-12-09 / Chase 117 / Bryce 126 / Zach 116),
# mirroring the date-column formatting used by the existing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A51").Value = 46000
$ws.Range("A51").NumberFormat = $ws.Range("A50").NumberFormat
$ws.Range("B51").Value = 117
$ws.Range("C51").Value = 126
$ws.Range("D51").Value = 116
